$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "maskin"
$ws.Range("B4").Value = 1500

$ws.Range("A5").Value = "kappa"
$ws.Range("B5").Value = 125

$ws.Range("A6").Value = "fil"
$ws.Range("B6").Value = 400

$ws.Range("A7").Value = "hätta"
$ws.Range("B7").Value = 124

$ws.Range("C7").Select()
